# Add a "Save" column (H) to the s_vals sheet:
#  - H1: header label "Save", styled like the other header cells (G1)
#  - H2: data value 1 (plain, unstyled, numeric)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: copy G1's formatting (bold, centered, bordered) onto H1,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Data cell: plain numeric value, no special formatting (matches B2:G2).
$ws.Range("H2").Value = 1
